$d = $word.ActiveDocument

# --- 1. Remove stray leading space before "Burnaby Mountain" (Team 1 Captains table, row 2) ---
$t1 = $d.Tables.Item(1)
$burnabyCell = $t1.Cell(2, 1)
$burnabyRange = $burnabyCell.Range
$leadingSpace1 = $d.Range($burnabyRange.Start, $burnabyRange.Start + 1)
if ($leadingSpace1.Text -eq " ") {
    $leadingSpace1.Delete()
}

# --- 2. Remove stray leading space before "Chilliwack" (Team 1 Captains table, row 3) ---
$t1b = $d.Tables.Item(1)
$chilliwackCell = $t1b.Cell(3, 1)
$chilliwackRange = $chilliwackCell.Range
$leadingSpace2 = $d.Range($chilliwackRange.Start, $chilliwackRange.Start + 1)
if ($leadingSpace2.Text -eq " ") {
    $leadingSpace2.Delete()
}

# --- 3 & 4. Move the "_GoBack" bookmark from after the page-break paragraph into the
# middle of the "604-536-6204" phone number (splitting that run into "60" / "4-536-6204") ---
$t2 = $d.Tables.Item(2)
$phoneCell = $t2.Cell(19, 3)
$phoneRange = $phoneCell.Range
$splitPos = $phoneRange.Start + 2
$bookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- 5. Center the paragraph containing the "marinedrive.bizwomen@gmail.com" hyperlink ---
$emailCell = $t2.Cell(20, 4)
$emailRange = $emailCell.Range
$emailRange.Paragraphs.Item(1).Alignment = 1
